$d = $word.ActiveDocument
$d.Content.Find.Execute("FROM ```Ecommerce`.firstCompany.orders", $true, $false, $false, $false, $false, $true, 1, $false, "FROM  ```Ecommerce`.firstCompany.orders", 2)
